$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column AM (culture_collection) - shifts cell/shared-string data left
$ws.Range("AM:AM").Delete()

# The comment boxes do not auto-shift with the column delete in this runtime,
# so re-stamp each remaining comment (AM15..CV15) with the text that now belongs
# there, then drop the trailing comment that lost its column (old CW15).
$commentUpdates = @(
    @('AM15', 'concentration of diether lipids; can include multiple types of diether lipids'),
    @('AN15', 'concentration of dissolved carbon dioxide'),
    @('AO15', 'concentration of dissolved hydrogen'),
    @('AP15', 'dissolved inorganic carbon concentration'),
    @('AQ15', 'concentration of dissolved organic carbon'),
    @('AR15', 'dissolved organic nitrogen concentration measured as; total dissolved nitrogen - NH4 - NO3 - NO2'),
    @('AS15', 'concentration of dissolved oxygen'),
    @('AT15', 'Plasmids that have significance phenotypic consequence'),
    @('AU15', 'measurement of glucosidase activity'),
    @('AV15', 'Health or disease status of sample at time of collection'),
    @('AW15', 'The natural (as opposed to laboratory) host to the organism from which the sample was obtained. Use the full taxonomic name, eg, "Homo sapiens".'),
    @('AX15', 'NCBI taxonomy ID of the host, e.g. 9606'),
    @('AY15', 'Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.'),
    @('AZ15', 'A locus tag prefix required for an annotated genome, http://www.ddbj.nig.ac.jp/sub/locus_tag-e.html'),
    @('BA15', 'concentration of magnesium'),
    @('BB15', 'measurement of mean friction velocity'),
    @('BC15', 'measurement of mean peak friction velocity'),
    @('BD15', 'methane (gas) amount or concentration at the time of sampling'),
    @('BE15', 'any other measurement performed or parameter collected, that is not listed here'),
    @('BF15', 'concentration of n-alkanes; can include multiple n-alkanes'),
    @('BG15', 'concentration of nitrate'),
    @('BH15', 'concentration of nitrite'),
    @('BI15', 'concentration of nitrogen (total)'),
    @('BJ15', 'concentration of organic carbon'),
    @('BK15', 'concentration of organic matter'),
    @('BL15', 'concentration of organic nitrogen'),
    @('BM15', 'total count of any organism per gram or volume of sample, should include name of organism followed by count; can include multiple organism counts'),
    @('BN15', 'oxygenation status of sample'),
    @('BO15', 'concentration of particulate organic carbon'),
    @('BP15', 'To what is the entity pathogenic'),
    @('BQ15', 'type of perturbation, e.g. chemical administration, physical disturbance, etc., coupled with time that perturbation occurred; can include multiple perturbation types'),
    @('BR15', 'concentration of petroleum hydrocarbon'),
    @('BS15', 'pH measurement'),
    @('BT15', 'concentration of phaeopigments; can include multiple phaeopigments'),
    @('BU15', 'concentration of phosphate'),
    @('BV15', 'concentration of phospholipid fatty acids; can include multiple values'),
    @('BW15', 'concentration of potassium'),
    @('BX15', 'pressure to which the sample is subject, in atmospheres'),
    @('BY15', 'redox potential, measured relative to a hydrogen cell, indicating oxidation or reduction potential'),
    @('BZ15', 'Primary publication or genome report in the form of pubmed ID, DOI or URL'),
    @('CA15', 'salinity measurement'),
    @('CB15', 'Method or device employed for collecting sample'),
    @('CC15', 'Processing applied to the sample during or after isolation'),
    @('CD15', 'Amount or size of sample (volume, mass or area) that was collected'),
    @('CE15', 'duration for which sample was stored'),
    @('CF15', 'location at which sample was stored, usually name of a specific freezer/room'),
    @('CG15', 'temperature at which sample was stored, e.g. -80'),
    @('CH15', 'volume (mL) or weight (g) of sample processed for DNA extraction'),
    @('CI15', 'concentration of silicate'),
    @('CJ15', 'sodium concentration'),
    @('CK15', 'unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.'),
    @('CL15', 'Identifier for the physical specimen. Use format: "[<institution-code>:[<collection-code>:]]<specimen_id>", eg, "UAM:Mamm:52179". Intended as a reference to the physical specimen that remains after it was analyzed. If the specimen was destroyed in the process of analysis, electronic images (e-vouchers) are an adequate substitute for a physical voucher specimen. Ideally the specimens will be deposited in a curated museum, herbarium, or frozen tissue collection, but often they will remain in a personal or laboratory collection for some time before they are deposited in a curated collection. There are three forms of specimen_voucher qualifiers. If the text of the qualifier includes one or more colons it is a ''structured voucher''. Structured vouchers include institution-codes (and optional collection-codes) taken from a controlled vocabulary maintained by the INSDC that denotes the museum or herbarium collection where the specimen resides, please visit the INSDC website, http://www.insdc.org/controlled-vocabulary-specimenvoucher-qualifier'),
    @('CM15', 'Information about the genetic distinctness of the lineage (eg., biovar, serovar)'),
    @('CN15', 'concentration of sulfate'),
    @('CO15', 'concentration of sulfide'),
    @('CP15', 'temperature of the sample at time of sampling'),
    @('CQ15', 'total carbon content'),
    @('CR15', 'total nitrogen content of the sample'),
    @('CS15', 'Definition for soil: total organic C content of the soil units of g C/kg soil. Definition otherwise: total organic carbon content'),
    @('CT15', 'Feeding position in food chain (eg., chemolithotroph)'),
    @('CU15', 'turbidity measurement'),
    @('CV15', 'water content measurement'),
)

foreach ($pair in $commentUpdates) {
    $ref = $pair[0]
    $newText = $pair[1]
    [void]$ws.Range($ref).Comment.Text($newText)
}

[void]$ws.Range("CW15").Comment.Delete()
